$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.883.36"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "3.854.45"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "598.44"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "168.62"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").Value = "0.666"
$ws.Range("E7").Value = "  -2.77%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "0.743"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  +4.01%  "
$ws.Range("D11").Value = "53.27"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").Value = "0.0000321"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").Value = "4.459.27"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").Value = "21.18"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").Value = "3.855.61"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "13.88"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("E18").Value = "  -4.92%  "
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "70.712.12"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "438.64"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "94.53"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("E24").Value = "  -4.90%  "
$ws.Range("D25").Value = "13.81"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("D26").Value = "11.52"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("E27").Value = "  -8.25%  "
$ws.Range("D28").Value = "5.97"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "10.48"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("E30").Value = "  +8.57%  "
$ws.Range("D31").Value = "34.98"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").Value = "13.51"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("D33").Value = "48.44"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  -4.56%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "68.76"
$ws.Range("E35").Value = "  -2.80%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0980"
$ws.Range("E36").Value = "  +8.56%  "
$ws.Range("D37").Value = "639.26"
$ws.Range("E37").Value = "  -4.59%  "
$ws.Range("D38").Value = "0.433"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "3.23"
$ws.Range("E42").Value = "  -4.59%  "
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  +8.09%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0470"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "3.11"
$ws.Range("E45").Value = "  +14.41%  "
$ws.Range("D46").Value = "10.17"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.91"
$ws.Range("E47").Value = "  -13.48%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.144"
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").Value = "2.945.08"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "3.32"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("D51").Value = "0.000276"
$ws.Range("E51").Value = "  +2.63%  "
